# Trade #105 (HighProbConvergence, DOWN) opened at 2026-02-18 00:19:19,
# and trade #76 (momentum, DOWN) closed via early_exit.
#
# This touches:
#   - Summary                : total trades / win rate roll-ups
#   - Strategy Status        : momentum strategy row roll-ups
#   - All Trades             : row 77 (closes trade #76) + new row 106 (trade #105)
#   - momentum                : row 11 (closes trade #76)
#   - HighProbConvergence     : new row 13 (trade #105)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 75        # Total Trades
$summary.Range("B9").Value = 49.33     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - momentum row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D11").Value = 9         # Trades
$status.Range("G11").Value = 11.11     # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - close out trade #76 (row 77) and append trade #105
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G77").Value = 0.01
$allTrades.Range("H77").Value = "CLOSED"
$allTrades.Range("K77").Value = 99.68000000000001
$allTrades.Range("L77").Value = "early_exit"
$allTrades.Range("M77").Value = 0.17

# Seed row 106 from row 105 (same shape: HighProbConvergence/DOWN, still OPEN)
# so date/time text cells keep their text type instead of Excel re-parsing
# "2026-02-18" into a real date when assigned directly.
$allTrades.Range("A105:Q105").Copy($allTrades.Range("A106:Q106"))
$allTrades.Range("A106").Value = 105
$allTrades.Range("C106").Value = "00:19:19"
$allTrades.Range("F106").Value = 0.01
$allTrades.Range("Q106").Value = "Mean reversion DOWN: price 3.55% above mean (z=3.00)"

# ---------------------------------------------------------------------------
# momentum sheet - close out trade #76 (row 11)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("G11").Value = 0.01
$momentum.Range("H11").Value = "CLOSED"
$momentum.Range("K11").Value = 99.68000000000001
$momentum.Range("P11").Value = "early_exit"
$momentum.Range("Q11").Value = 0.17

# ---------------------------------------------------------------------------
# HighProbConvergence sheet - append trade #105 (row 13), seeded from row 12
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")

$hpc.Range("A12:Q12").Copy($hpc.Range("A13:Q13"))
$hpc.Range("A13").Value = 105
$hpc.Range("C13").Value = "00:19:19"
$hpc.Range("F13").Value = 0.01
$hpc.Range("O13").Value = "Mean reversion DOWN: price 3.55% above mean (z=3.00)"
